$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.855280357757138
$ws.Range("C2").Value = 0.8836601307189543
$ws.Range("D2").Value = 0.868814192343604

$ws.Range("B3").Value = 0.8261327561327562
$ws.Range("C3").Value = 0.8489473684210527
$ws.Range("D3").Value = 0.8359820282413348

$ws.Range("B4").Value = 0.8305627705627707
$ws.Range("C4").Value = 0.9263157894736842
$ws.Range("D4").Value = 0.8755659787367105

$ws.Range("B5").Value = 0.8241946778711483
$ws.Range("C5").Value = 0.7970588235294118
$ws.Range("D5").Value = 0.8083836467138175

$ws.Range("B6").Value = 0.9215686274509803
$ws.Range("C6").Value = 0.7852380952380952
$ws.Range("D6").Value = 0.8471605208447313
